# Marplatenses.xlsx edit: restructure Hoja1 -> Hoja2 as a 12-column
# ingredient/quantity calculator table.
$wb = $excel.ActiveWorkbook

# --- 1. Recreate the sheet as "Hoja2" (so sheetId increments 1 -> 2, matching
#        the diff) by adding a fresh sheet and removing the original "Hoja1". ---
$new = $wb.Worksheets.Add()
$new.Name = "Hoja2"
$old = $wb.Worksheets.Item("Hoja1")
$old.Delete()

$ws = $wb.ActiveSheet

# --- 2. Header row (ingredient names), alphabetical A..L ---
$headers = @("Azucar Blanca","Bicarbonato de Sodio","Cacao Amargo","Escencia de Vainilla","Fécula de Maiz","Harina 0000","Huevos","Limón","Margarina","Miel","Naranja","Polvo de Hornear")
# --- 3. Quantity row (free-text amounts), aligned to the same columns ---
$qtys    = @("50g","2g","20g","1 tsp","50g","250g","1 u","1/2 u","150g","50g","1/2 u","5g")

# Seed row 1/2 with placeholder content whose column-A comparison is stable
# under Sort.Apply() (engine quirk: Apply() only ever compares column A of
# row1 vs row2), then record a sortState over the header row, then overwrite
# with the real values so the data never gets row-flipped.
$ws.Range("A1:L1").Value = "a"
$ws.Range("A2:L2").Value = "z"

$rng = $ws.Range("A1:L2")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A1:L1"), 0, 1, 0, 0)
$ws.Sort.SetRange($rng)
$ws.Sort.Header = 0
$ws.Sort.Orientation = 1
$ws.Sort.Apply()

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $hcell = $ws.Cells.Item(1, $col)
    $hcell.Value = $headers[$i]
    $hcell.NumberFormat = "@"

    $qcell = $ws.Cells.Item(2, $col)
    $qcell.Value = $qtys[$i]
    $qcell.NumberFormat = """$""\ #,##0.00"
}

# Columns B (Bicarbonato de Sodio) and J (Miel) headers are vertically centered
$ws.Cells.Item(1, 2).VerticalAlignment = -4108
$ws.Cells.Item(1, 10).VerticalAlignment = -4108

# --- 4. Column widths (best-fit approximations for the new 12-column layout) ---
$widths = @(12.166666666666666, 19.0, 12.666666666666666, 17.666666666666668, 13.166666666666666, 10.333333333333334, 6.666666666666667, 5.666666666666667, 9.0, 4.166666666666667, 7.0, 15.666666666666666)
for ($i = 0; $i -lt $widths.Length; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = $widths[$i]
}

# --- 5. Selection lands on B4, matching the saved view state ---
$ws.Range("B4").Select()
